$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- Sheet: 展览 ----
$ws1.Range("F2").Value = 336
$ws1.Range("C3").Value = "萍乡·BM次元盛典运动番only（取消）"
$ws1.Range("F3").Value = 279
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F4").Value = 1233
$ws1.Range("C5").Value = "宜春·BM次元盛典运动番only（取消）"
$ws1.Range("G5").Value = "不可售"
$ws1.Range("C7").Value = "鹰潭·BM次元盛典运动番only（取消）"
$ws1.Range("G7").Value = "不可售"
$ws1.Range("C8").Value = "赣州·BM次元盛典运动番only（取消）"
$ws1.Range("G8").Value = "不可售"
$ws1.Range("F9").Value = 142
$ws1.Range("F10").Value = 3437
$ws1.Range("F11").Value = 124
$ws1.Range("F12").Value = 84
$ws1.Range("F14").Value = 39
$ws1.Range("F15").Value = 54
$ws1.Range("F16").Value = 589
$ws1.Range("F17").Value = 81
$ws1.Range("F18").Value = 714
$ws1.Range("F20").Value = 118
$ws1.Range("F23").Value = 64
$ws1.Range("F24").Value = 2559
$ws1.Range("F25").Value = 5066
$ws1.Range("F28").Value = 476
$ws1.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202406/meKBC0hU1719222126375.jpeg"
$ws1.Range("F29").Value = 1295
$ws1.Range("F30").Value = 280
$ws1.Range("F31").Value = 2229
$ws1.Range("F34").Value = 79
$ws1.Range("F35").Value = 103
$ws1.Range("F36").Value = 170
$ws1.Range("F37").Value = 309
$ws1.Range("F38").Value = 458
$ws1.Range("F39").Value = 789
$ws1.Range("F40").Value = 29
$ws1.Range("F42").Value = 34
$ws1.Range("F43").Value = 470

# ---- Sheet: 全部类型 ----
$ws4.Range("F2").Value = 336
$ws4.Range("C3").Value = "萍乡·BM次元盛典运动番only（取消）"
$ws4.Range("F3").Value = 279
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F4").Value = 1234
$ws4.Range("C5").Value = "宜春·BM次元盛典运动番only（取消）"
$ws4.Range("G5").Value = "不可售"
$ws4.Range("C7").Value = "鹰潭·BM次元盛典运动番only（取消）"
$ws4.Range("G7").Value = "不可售"
$ws4.Range("C8").Value = "赣州·BM次元盛典运动番only（取消）"
$ws4.Range("G8").Value = "不可售"
$ws4.Range("F9").Value = 142
$ws4.Range("F10").Value = 3438
$ws4.Range("F11").Value = 124
$ws4.Range("F12").Value = 84
$ws4.Range("F15").Value = 39
$ws4.Range("F16").Value = 54
$ws4.Range("F17").Value = 589
$ws4.Range("F18").Value = 81
$ws4.Range("F19").Value = 714
$ws4.Range("F21").Value = 118
$ws4.Range("F24").Value = 64
$ws4.Range("F25").Value = 2559
$ws4.Range("F26").Value = 5066
$ws4.Range("F29").Value = 476
$ws4.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202406/meKBC0hU1719222126375.jpeg"
$ws4.Range("F30").Value = 1295
$ws4.Range("F31").Value = 280
$ws4.Range("F32").Value = 2229
$ws4.Range("F35").Value = 79
$ws4.Range("F36").Value = 103
$ws4.Range("F37").Value = 170
$ws4.Range("F38").Value = 309
$ws4.Range("F39").Value = 458
$ws4.Range("F40").Value = 789
$ws4.Range("F41").Value = 29
$ws4.Range("F43").Value = 34
$ws4.Range("F44").Value = 470
